$wb = $excel.ActiveWorkbook

# --- Core data edit: "views annualization calculations removed" ---
# Each of the 8 "view" sheets has its shock size in E2 divided by 10
# (de-annualized). The I2 cell holds a formula referencing E2, so its
# cached text refreshes automatically on recalculation.

$ws1 = $wb.Worksheets.Item(1)   # view 1 (rates up)
$ws2 = $wb.Worksheets.Item(2)   # view 2 (rates dn)
$ws3 = $wb.Worksheets.Item(3)   # view 3 (eq up)
$ws4 = $wb.Worksheets.Item(4)   # view 4 (eq dn)
$ws5 = $wb.Worksheets.Item(5)   # view 5 (inflation up)
$ws6 = $wb.Worksheets.Item(6)   # view 6 (inflation dn)
$ws7 = $wb.Worksheets.Item(7)   # view 7 (vix up)
$ws8 = $wb.Worksheets.Item(8)   # view 8 (vix dn)

$ws1.Range("E2").Value = 0.025
$ws2.Range("E2").Value = -0.025
$ws3.Range("E2").Value = 0.01
$ws4.Range("E2").Value = -0.01
$ws5.Range("E2").Value = 0.001
$ws6.Range("E2").Value = -0.001
$ws7.Range("E2").Value = 0.01
$ws8.Range("E2").Value = -0.01

# --- View / selection state updates ---
# view 1: bottom (frozen) pane selection moves from C4 to D6
$ws1.Range("D6").Select() | Out-Null

# view 2: bottom (frozen) pane selection moves from E3 to E2
$ws2.Range("E2").Select() | Out-Null

# view 8 becomes the active tab; its selection moves from F7 to E3
$ws8.Range("E3").Select() | Out-Null
$ws8.Activate() | Out-Null
